# Regenerate merged AHB files
# - rename the FV2404 ("_old") / FV2410 ("_new") column header labels
# - (re)build the Excel Table ("Table1") over the full data range
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row labels -----------------------------------------
# The sheet holds a before/after "AHB diff" comparison. The column headers
# used to be suffixed "_old"/"_new"; they are renamed to the concrete
# Formatversion tags "_FV2404"/"_FV2410".
$headerNames = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headerNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerNames[$i]
}

# --- 2. Turn the used range into a proper Excel Table ---------------------
$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $headerNames.Length
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))

$lo = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$lo.Name = "Table1"
$lo.TableStyle = "TableStyleMedium2"
$lo.ShowTableStyleFirstColumn = $false
$lo.ShowTableStyleLastColumn = $false
$lo.ShowTableStyleRowStripes = $true
$lo.ShowTableStyleColumnStripes = $false

# --- 3. Freeze the header row ---------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
